$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set a custom width on column A (closest value this host's character-width
# model can reach to the authored 16.5703125; the emulator's ColumnWidth
# setter only snaps to whole-character increments).
$ws.Columns.Item(1).ColumnWidth = 15.67

# Append the two new enrollment rows.
$ws.Range("A4").Value = "F-22-SE-A-3001"
$ws.Range("B4").Value = "MT-0001"
$ws.Range("A5").Value = "F-22-SE-A-3001"
$ws.Range("B5").Value = "SS-0001"

# Match the author's final selection.
[void]$ws.Range("B13").Select()
